$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 4 (SK35SMA diode) ----
$ws.Range("A4").Value = "D?"
$ws.Range("C4").Value = "SK35SMA"
$ws.Range("D4").Value = "SMA"
$ws.Range("F4").Value = "Dioda: prostownicza Schottky; SMD; 50V; 3A; SMA"
$ws.Range("J4").Value = "https://www.tme.eu/pl/details/sk35sma-dio/diody-schottky-smd/diotec-semiconductor/sk35sma/"
$ws.Range("B4").Value = "DIOTEC SEMICONDUCTOR"
$ws.Range("E4").Value = "SMT"
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 0.3496
$ws.Range("I4").Formula = "=G4*H4"

# ---- Row 5 (A4403GEUTR-T regulator) ----
$ws.Range("A5").Value = "U?"
$ws.Range("B5").Value = "ALLEGRO"
$ws.Range("C5").Value = "A4403GEUTR-T"
$ws.Range("D5").Value = "QFN16"
$ws.Range("F5").Value = "PMIC; przetwornica DC/DC; Upracy: 9÷46V; Uwyj: 46V; QFN16; buck"
$ws.Range("J5").Value = "https://www.tme.eu/pl/details/a4403geutr-t/regulatory-napiecia-uklady-dc-dc/allegro-microsystems/"
$ws.Range("E5").Value = "SMT"
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 6.68
$ws.Range("I5").Formula = "=G5*H5"

# ---- Row 6 (DJNR6045-6R3-S inductor) ----
$ws.Range("C6").Value = "DJNR6045-6R3-S"
$ws.Range("A6").Value = "L?"
$ws.Range("B6").Value = "FERROCORE"
$ws.Range("F6").Value = "Dławik: drutowy; SMD; 6,3uH; Ipracy: 3A; 36mΩ; 6x6x4,5mm; ±20%"
$ws.Range("J6").Value = "https://www.tme.eu/pl/details/djnr6045-6r3-s/dlawiki-smd-mocy/ferrocore/"
$ws.Range("E6").Value = "SMT"
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 0.7064
$ws.Range("I6").Formula = "=G6*H6"

# ---- Row 7 (TMC2209-LA-T stepper driver) ----
$ws.Range("C7").Value = "TMC2209-LA-T"
$ws.Range("A7").Value = "U?"
$ws.Range("B7").Value = "TRINAMIC"
$ws.Range("D7").Value = "QFN28"
$ws.Range("F7").Value = "Sterownik/kontroler silnika krokowego, 4.75V do 29V, 2A/1 wyjście, QFN-28"
$ws.Range("J7").Value = "https://pl.farnell.com/trinamic/tmc2209-la-t/motor-driver-stepper-qfn-28/dp/3131535?st=tmc2209"
$ws.Range("E7").Value = "SMT"
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 12.36
$ws.Range("I7").Formula = "=G7*H7"

# ---- Row 25 (Total) - entered before filling in the rest of rows 8 & 9 ----
$ws.Range("G25").Value = "Razem"
$ws.Range("I25").Formula = "=SUM(I2:I24)"

# ---- Row 8 (DJNR6045-150 inductor) ----
$ws.Range("A8").Value = "L?"
$ws.Range("B8").Value = "FERROCORE"
$ws.Range("C8").Value = "DJNR6045-150"
$ws.Range("F8").Value = "Dławik: drutowy; SMD; 15uH; Ipracy: 2,3A; 100,1mΩ; 6x5,9x4,5mm; ±20%"
$ws.Range("J8").Value = "https://www.tme.eu/pl/details/djnr6045-150/dlawiki-smd-mocy/ferrocore/"
$ws.Range("E8").Value = "SMT"
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 0.6992
$ws.Range("I8").Formula = "=G8*H8"

# ---- Row 9 (MCP16301T-I/CHY regulator) ----
$ws.Range("C9").Value = "MCP16301T-I/CHY"
$ws.Range("B9").Value = "MICROCHIP TECHNOLOGY"
$ws.Range("D9").Value = "SOT23-6"
$ws.Range("F9").Value = "PMIC; przetwornica DC/DC; Upracy: 4÷30V; Uwyj: 2÷15V; SOT23-6; buck"
$ws.Range("J9").Value = "https://www.tme.eu/pl/details/mcp16301t-i_chy/regulatory-napiecia-uklady-dc-dc/microchip-technology/"
$ws.Range("A9").Value = "U?"
$ws.Range("E9").Value = "SMT"
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 3.84
$ws.Range("I9").Formula = "=G9*H9"

# Update sheet view selection
$ws.Range("F14").Select()
